$wb = $excel.ActiveWorkbook

# "Metadata" sheet holds the properties table (Property in col A, Value in col B)
$meta = $wb.Worksheets.Item("Metadata")

# Update the Date value (row 8, col B) to the new timestamp
$meta.Range("B8").Value = "2025-07-17T14:35:50+00:00"

# Fill in the previously empty "Title" value (row 5, col B) with "Competence"
$meta.Range("B5").Value = "Competence"
